$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently holds 4 data rows (rows 2-5). A new data row needs to be
# inserted as the new first data row (row 2), pushing the existing rows down
# by one (old row 2 -> 3, old row 3 -> 4, old row 4 -> 5, old row 5 -> 6).
# Shift the data down (starting from the bottom so we don't overwrite
# anything before it has been copied down).
$ws.Range("A6:T6").Value2 = $ws.Range("A5:T5").Value2
$ws.Range("A5:T5").Value2 = $ws.Range("A4:T4").Value2
$ws.Range("A4:T4").Value2 = $ws.Range("A3:T3").Value2
$ws.Range("A3:T3").Value2 = $ws.Range("A2:T2").Value2

# Row 6 is brand new, so copy the date number format used in column D from the
# row above it (all data rows share the same date format in column D).
$ws.Range("D6").NumberFormat = $ws.Range("D5").NumberFormat

# Fill in the new first data row (row 2) with the new record's values.
$ws.Range("A2").Value2 = 5
$ws.Range("B2").Value2 = "Macroferia Regional de Talca"
$ws.Range("C2").Value2 = "Maule"
$ws.Range("D2").Value2 = 45250
$ws.Range("E2").Value2 = 7
$ws.Range("F2").Value2 = "Fruta"
$ws.Range("G2").Value2 = 100104
$ws.Range("H2").Value2 = "Frutos de pepita"
$ws.Range("I2").Value2 = 100104004
$ws.Range("J2").Value2 = "Níspero"
$ws.Range("K2").Value2 = "Golden Nugget"
$ws.Range("L2").Value2 = "Primera"
$ws.Range("M2").Value2 = 100
$ws.Range("N2").Value2 = 20000
$ws.Range("O2").Value2 = 20000
$ws.Range("P2").Value2 = 20000
$ws.Range("Q2").Value2 = "`$/bandeja 10 kilos"
$ws.Range("R2").Value2 = "Provincia de Limarí"
$ws.Range("S2").Value2 = 2000
$ws.Range("T2").Value2 = 10
